$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(992).Insert()

$ws.Range("A992").Value = 3
$ws.Range("B992").Value = "Femacal de La Calera"
$ws.Range("C992").Value = "Coquimbo"
$ws.Range("D992").Value = 45212
$ws.Range("E992").Value = 5
$ws.Range("F992").Value = 100112006
$ws.Range("G992").Value = "Repollo"
$ws.Range("H992").Value = "Crespo record"
$ws.Range("I992").Value = "Primera"
$ws.Range("J992").Value = 2700
$ws.Range("K992").Value = 700
$ws.Range("L992").Value = 800
$ws.Range("M992").Value = 756
$ws.Range("N992").Value = "$/unidad"
$ws.Range("O992").Value = "Provincia de Quillota"
$ws.Range("P992").Value = 756
$ws.Range("Q992").Value = 1
$ws.Range("R992").Value = "Hortaliza"
